# Apply the "Sara Caballeria" attitude data update on the "Ninja" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ninja")
$ws.Activate()

# Correct the student name in A5 (adds a new shared string rather than
# editing the existing "Sara Carballeira" entry).
$ws.Range("A5").Value = "Sara Caballeria"

# Fill in the new T7 (G) / T8 (H) attitude scores for each student.
$ws.Range("G2").Value = 0.5
$ws.Range("H2").Value = 0.5

$ws.Range("G5").Value = 1
$ws.Range("H5").Value = 1

$ws.Range("G7").Value = 0.75
$ws.Range("H7").Value = 0.75

$ws.Range("G8").Value = 0.75
$ws.Range("H8").Value = 0.75

$ws.Range("G9").Value = 1
$ws.Range("H9").Value = 1

$ws.Range("G11").Value = 1
$ws.Range("H11").Value = 1

$ws.Range("G14").Value = 1
$ws.Range("H14").Value = 1

$ws.Range("G16").Value = 0
$ws.Range("H16").Value = 0

$ws.Range("G17").Value = 1
$ws.Range("H17").Value = 1

$ws.Range("G18").Value = 1
$ws.Range("H18").Value = 1

$ws.Range("G19").Value = 0.25
$ws.Range("H19").Value = 0.25

$ws.Range("G20").Value = 1
$ws.Range("H20").Value = 1

$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 1

$ws.Range("G23").Value = 1
$ws.Range("H23").Value = 1

$ws.Range("G24").Value = 0.5
$ws.Range("H24").Value = 0.5

$ws.Range("G25").Value = 0.5
$ws.Range("H25").Value = 0.5

$ws.Range("G26").Value = 0.75
$ws.Range("H26").Value = 0.75

$ws.Range("G27").Value = 0.75
$ws.Range("H27").Value = 0.75

$ws.Range("G28").Value = 0.75
$ws.Range("H28").Value = 0.75

# Leave the active selection where the author left it after editing.
$ws.Range("P23").Select() | Out-Null
